# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G (header "K") values are re-written with new computed values
# (strikeouts instead of raw strike-count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 8
    6  = 4
    7  = 6
    8  = 10
    9  = 7
    10 = 5
    11 = 4
    12 = 4
    13 = 4
    14 = 4
    15 = 5
    16 = 6
    17 = 9
    18 = 4
    19 = 7
    20 = 6
    21 = 11
    22 = 8
    23 = 9
    24 = 9
    25 = 5
    26 = 8
    27 = 1
    28 = 3
    29 = 2
    30 = 5
    31 = 3
    32 = 4
    33 = 3
    34 = 3
    35 = 6
    36 = 4
    37 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
